$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Week 13 (column N) day-after inputs: fill in actual results that were
# previously marked as "A" (absent) pending the coach's update.
$ws.Range("N3").Value = "W"
$ws.Range("N4").Value = "L"
$ws.Range("N5").Value = "L"
$ws.Range("N6").Value = "NA"
$ws.Range("N7").Value = "DNP"
$ws.Range("N8").Value = "DNP"
$ws.Range("N9").Value = "NA"
$ws.Range("N10").Value = "L"

$ws.Range("N15").Value = "NA"
$ws.Range("M16").Value = "L"
$ws.Range("N16").Value = "W"
$ws.Range("N17").Value = "W"
$ws.Range("N18").Value = "W"
$ws.Range("N19").Value = "L"
$ws.Range("N20").Value = "DNP"
$ws.Range("N21").Value = "NA"
$ws.Range("N22").Value = "W"

# Update the last active selection to reflect where entry ended.
$ws.Range("S19").Select()
